$d = $word.ActiveDocument

# --- 1. Remove the existing "_GoBack" bookmark from the "Planned activities:"
#        paragraph; it will be recreated at the end of the new last paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Insert the three new paragraphs right after "Highlights:".
$highlights = $d.Paragraphs(3)
$r = $highlights.Range
$r.Collapse(0)              # wdCollapseEnd
$r.InsertParagraphAfter()

$p4 = $d.Paragraphs(4)
$p4.Range.InsertBefore("Georgi worked on Renting GUI")

$r2 = $p4.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()

$p5 = $d.Paragraphs(5)
$p5.Range.InsertBefore("Angel worked on issues with the tickets at the website")

$r3 = $p5.Range
$r3.Collapse(0)
$r3.InsertParagraphAfter()

$p6 = $d.Paragraphs(6)
$p6.Range.InsertBefore("Mikaeil practiced C#")

# --- 3. Re-add the "_GoBack" bookmark at the end of the new last paragraph's
#        text (zero-length, right before the paragraph mark) -- same relative
#        position it held before the edit.
$p6b = $d.Paragraphs(6)
$textEnd = $p6b.Range.End - 1   # position right after "C#", before the pilcrow

# Work around a zero-length-range bookmark placement quirk: insert a
# throw-away placeholder character, bookmark the (non-empty) range that
# contains it, then delete the placeholder. The bookmark collapses to the
# correct empty position once its contents are removed.
$ins = $d.Range($textEnd, $textEnd)
$ins.InsertAfter("X")

$bmRange = $d.Range($textEnd, $textEnd + 1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$placeholder = $d.Range($textEnd, $textEnd + 1)
$placeholder.Delete() | Out-Null
